$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")

# Rename archive columns - append "(AH)" suffix to the model names in row 1
$ws.Range("B1").Value = "ESP(AH)"
$ws.Range("C1").Value = "BLAST(AH)"
$ws.Range("D1").Value = "DOE2(AH)"
$ws.Range("E1").Value = "SRES/SUN(AH)"
$ws.Range("F1").Value = "SERIRES(AH)"
$ws.Range("G1").Value = "S3PAS(AH)"
$ws.Range("H1").Value = "TRNSYS(AH)"
$ws.Range("I1").Value = "TASE(AH)"

# Minor fix of distance_% values for CASE 960 and CASE 210
$ws.Range("M10").Value = 2.9
$ws.Range("M13").Value = 6.1
